$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 141 (shifts existing rows 141..256 down to 142..257)
$ws.Rows.Item(141).Insert()

# Populate the newly inserted row 141 with the new data record
$ws.Range("A141").Value = 3
$ws.Range("B141").Value = "Femacal de La Calera"
$ws.Range("C141").Value = "Coquimbo"
$ws.Range("D141").Value = 44483
$ws.Range("E141").Value = 5
$ws.Range("F141").Value = 100112028
$ws.Range("G141").Value = "Sandia"
$ws.Range("H141").Value = "Sin especificar"
$ws.Range("I141").Value = "Primera"
$ws.Range("J141").Value = 320
$ws.Range("K141").Value = 800
$ws.Range("L141").Value = 800
$ws.Range("M141").Value = 800
$ws.Range("N141").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O141").Value = "Perú"
$ws.Range("P141").Value = 800
$ws.Range("Q141").Value = 1
$ws.Range("R141").Value = "Hortaliza"
